$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New booking data for rows 2-5
$data = @(
    @(476042, "Car Assigned",    45956.12550925926, "Weekly", "Old", 1013.32, 0.08, 85.1352,  "D95",  0.000325, 32.93),
    @(477148, "Car Assigned",    45956.12550925926, "Daily",  "Old", 106.18,  0.1,  11.798,   "D95",  0.0003,   3.19),
    @(477159, "Vendor Assigned", 45956.12550925926, "Daily",  "Old", 101.67,  0.1,  11.298,   "D100", 0.0003,   3.05),
    @(477252, "Vendor Assigned", 45956.12550925926, "Daily",  "Old", 690.9,   0.1,  73.8924,  "D100", 0.0003,   20.73)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $vals[$j]
    }
}

# Remove the two old trailing booking rows (old rows 6 and 7); this shifts
# the old empty row 8 up to become row 6
$ws.Rows("6:7").Delete()
